# Refresh the cryptos price snapshot (GitHub Actions scrape update).
#
# Column D ("Price") holds plain text in the workbook (e.g. "3.273.60" uses
# dots as thousands separators, so it can never round-trip as a number). Several
# of the new prices (e.g. "0.999", "8.81") *do* look like plain numbers, so a
# leading apostrophe character is written in front of them (an embedded single
# quote, doubled per PowerShell single-quoted-string escaping rules) so Excel
# keeps storing/treating them as text, matching the original inline-string
# cells instead of silently converting them to a number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.395.09'
$ws.Range('E2').Value = '  +0.69%  '
$ws.Range('D3').Value = '3.273.60'
$ws.Range('E3').Value = '  +2.86%  '
$ws.Range('D4').Value = '''0.999'
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Value = '''609.61'
$ws.Range('E5').Value = '  +0.83%  '
$ws.Range('D6').Value = '''158.71'
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('D8').Value = '3.272.35'
$ws.Range('E8').Value = '  +2.80%  '
$ws.Range('E10').Value = '  +2.29%  '
$ws.Range('E11').Value = '  +1.56%  '
$ws.Range('D12').Value = '''0.497'
$ws.Range('E12').Value = '  -3.05%  '
$ws.Range('D13').Value = '''0.0000272'
$ws.Range('E13').Value = '  +3.01%  '
$ws.Range('D14').Value = '''38.86'
$ws.Range('E14').Value = '  +1.72%  '
$ws.Range('D15').Value = '3.809.16'
$ws.Range('E15').Value = '  +2.99%  '
$ws.Range('D16').Value = '66.428.56'
$ws.Range('E16').Value = '  +0.69%  '
$ws.Range('D17').Value = '3.273.53'
$ws.Range('E17').Value = '  +3.04%  '
$ws.Range('D18').Value = '''7.38'
$ws.Range('E18').Value = '  +0.01%  '
$ws.Range('E19').Value = '  +1.20%  '
$ws.Range('D20').Value = '''503.46'
$ws.Range('E20').Value = '  -0.85%  '
$ws.Range('D21').Value = '''15.44'
$ws.Range('E21').Value = '  +0.94%  '
$ws.Range('D22').Value = '''0.754'
$ws.Range('E22').Value = '  +3.53%  '
$ws.Range('D23').Value = '''8.13'
$ws.Range('E23').Value = '  +1.57%  '
$ws.Range('D24').Value = '''14.75'
$ws.Range('E24').Value = '  -0.12%  '
$ws.Range('D25').Value = '''86.57'
$ws.Range('E25').Value = '  +2.71%  '
$ws.Range('E26').Value = '  +0.03%  '
$ws.Range('D27').Value = '''3.03'
$ws.Range('E27').Value = '  +1.47%  '
$ws.Range('E28').Value = '  +0.26%  '
$ws.Range('D29').Value = '''0.137'
$ws.Range('E29').Value = '  +55.36%  '
$ws.Range('D30').Value = '''2.39'
$ws.Range('E30').Value = '  +0.62%  '
$ws.Range('E31').Value = '  -1.79%  '
$ws.Range('D32').Value = '''2.86'
$ws.Range('E32').Value = '  -4.78%  '
$ws.Range('D33').Value = '''27.98'
$ws.Range('E33').Value = '  +0.28%  '
$ws.Range('E34').Value = '  +0.05%  '
$ws.Range('E35').Value = '  -1.83%  '
$ws.Range('D36').Value = '''6.48'
$ws.Range('E36').Value = '  +0.03%  '
$ws.Range('D37').Value = '''3.52'
$ws.Range('E37').Value = '  +23.19%  '
$ws.Range('D38').Value = '0.0₃0799'
$ws.Range('E38').Value = '  +15.23%  '
$ws.Range('B39').Value = 'OKB'
$ws.Range('C39').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D39').Value = '''55.69'
$ws.Range('E39').Value = '  +0.92%  '
$ws.Range('B40').Value = 'Bittensor'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D40').Value = '''500.99'
$ws.Range('E40').Value = '  -2.48%  '
$ws.Range('E41').Value = '  +2.51%  '
$ws.Range('E42').Value = '  +2.38%  '
$ws.Range('D43').Value = '''8.81'
$ws.Range('D44').Value = '''2.60'
$ws.Range('E44').Value = '  +4.99%  '
$ws.Range('D45').Value = '3.007.62'
$ws.Range('E45').Value = '  +6.20%  '
$ws.Range('D46').Value = '''0.292'
$ws.Range('E46').Value = '  -2.01%  '
$ws.Range('D47').Value = '''28.96'
$ws.Range('E47').Value = '  +4.14%  '
$ws.Range('D48').Value = '''2.48'
$ws.Range('E48').Value = '  +4.78%  '
$ws.Range('E49').Value = '  +2.27%  '
$ws.Range('E50').Value = '  +0.00%  '
$ws.Range('D51').Value = '''121.30'
$ws.Range('E51').Value = '  -1.26%  '

Write-Host "Applied cryptos update"
